# Daily attendance processing - reorder the "Recorded By" (column G) names.
#
# For every data row, if the "Recorded By" cell lists more than one
# recorder (comma-separated) and the list does not already start with
# "System"/"system", reverse the order of the names so that the
# System/system entry ends up first (mirrors the author's manual fixup).
# Cells with a single recorder, or that already start with System, are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByCol = 7   # Column G: "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $text = [string]$cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ','
    if ($parts.Count -le 1) {
        continue
    }

    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    if ($trimmedParts[0].ToLower() -eq "system") {
        continue
    }

    $reversedIndices = ($trimmedParts.Count - 1)..0
    $reversedParts = $trimmedParts[$reversedIndices]
    $newText = [string]::Join(", ", $reversedParts)

    $cell.Value = $newText
}
